$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, and week-of dates) ---
$a8 = $ws.Range("A8")
$a8Text = $a8.Text
$a8.Characters($a8Text.Length, 1).Text = "3"

$c9 = $ws.Range("C9")
$c9.Characters(46, 9).Text = "1/22/2023"
$c9.Characters(27, 8).Text = "1/16/2023"

# --- Cells changing from numeric to text placeholders ("0" / "***.*") ---
function Set-TextFromTemplate($targetAddr, $templateAddr) {
    $ws.Range($templateAddr).Copy()
    $ws.Range($targetAddr).PasteSpecial(-4163)
    $ws.Range($templateAddr).Copy()
    $ws.Range($targetAddr).PasteSpecial(-4122)
}

Set-TextFromTemplate "C14" "D14"
Set-TextFromTemplate "D27" "C22"
Set-TextFromTemplate "E27" "E14"
Set-TextFromTemplate "C28" "D28"
Set-TextFromTemplate "C29" "D29"

# --- Cells changing from text placeholder to numeric ---
function Set-NumberFromTemplate($targetAddr, $templateAddr, $value) {
    $ws.Range($templateAddr).Copy()
    $ws.Range($targetAddr).PasteSpecial(-4122)
    $ws.Range($targetAddr).Value = $value
}

Set-NumberFromTemplate "L26" "K26" -100

# --- Plain numeric value updates ---
$ws.Range("G15").Value = 4
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 88.888888888888
$ws.Range("I16").Value = 11
$ws.Range("J16").Value = 8
$ws.Range("K16").Value = 37.5
$ws.Range("L16").Value = 266.666666666667
$ws.Range("M16").Value = -31.25
$ws.Range("N16").Value = -86.075949367088
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 27.777777777777
$ws.Range("I17").Value = 16
$ws.Range("J17").Value = 14
$ws.Range("K17").Value = 14.285714285714
$ws.Range("L17").Value = -11.111111111111
$ws.Range("M17").Value = 77.777777777777
$ws.Range("N17").Value = 23.076923076923
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -16.666666666666
$ws.Range("I18").Value = 14
$ws.Range("J18").Value = 18
$ws.Range("K18").Value = -22.222222222222
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -53.333333333333
$ws.Range("N18").Value = -90
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -80
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 58
$ws.Range("H19").Value = -48.275862068965
$ws.Range("I19").Value = 23
$ws.Range("J19").Value = 40
$ws.Range("K19").Value = -42.5
$ws.Range("L19").Value = -30.30303030303
$ws.Range("M19").Value = -34.285714285714
$ws.Range("N19").Value = -28.125
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 16.666666666666
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = 18.181818181818
$ws.Range("I20").Value = 17
$ws.Range("J20").Value = 18
$ws.Range("K20").Value = -5.555555555555
$ws.Range("L20").Value = 21.428571428571
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -91.70731707317
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = -30.76923076923
$ws.Range("F21").Value = 118
$ws.Range("G21").Value = 134
$ws.Range("H21").Value = -11.940298507462
$ws.Range("I21").Value = 82
$ws.Range("J21").Value = 99
$ws.Range("K21").Value = -17.171717171717
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -23.364485981308
$ws.Range("N21").Value = -82.627118644067
$ws.Range("F22").Value = 2
$ws.Range("M22").Value = -75
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = -16
$ws.Range("F24").Value = 102
$ws.Range("G24").Value = 94
$ws.Range("H24").Value = 8.510638297872
$ws.Range("I24").Value = 79
$ws.Range("J24").Value = 72
$ws.Range("K24").Value = 9.722222222222
$ws.Range("L24").Value = -26.851851851851
$ws.Range("M24").Value = 21.538461538461
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -33.333333333333
$ws.Range("I25").Value = 25
$ws.Range("J25").Value = 22
$ws.Range("K25").Value = 13.636363636363
$ws.Range("L25").Value = 8.695652173913
$ws.Range("M25").Value = -39.024390243902
$ws.Range("G26").Value = 4
$ws.Range("F27").Value = 7
$ws.Range("H27").Value = 75
$ws.Range("I27").Value = 5
$ws.Range("K27").Value = 25
